$d = $word.ActiveDocument

# Locate the "LOQ4084: Fenômenos de Transporte II (Requisito fraco)" paragraph,
# which is immediately followed by a blank paragraph, a "Ver no Jupiter..."
# paragraph and a "© 2020 ..." footer paragraph that must be removed.
$anchorText = "LOQ4084: Fenômenos de Transporte II (Requisito fraco)"

$searchRange = $d.Content
$found = $searchRange.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph text"
}

# Resolve the paragraph index that contains the found range.
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $searchRange.Start -and $p.Range.End -ge $searchRange.End) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq 0) {
    throw "Could not resolve anchor paragraph index"
}

# The three paragraphs right after the anchor are the ones to delete:
#   1) an empty paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
$firstToDelete = $d.Paragraphs.Item($anchorIndex + 1)
$lastToDelete = $d.Paragraphs.Item($anchorIndex + 3)

$deleteRange = $d.Range($firstToDelete.Range.Start, $lastToDelete.Range.End)
$deleteRange.Delete()
